# Aggiornamento fino a 28/06 incluso
# Appends daily rows (270-301) covering date serials 44344..44375
# (2021-05-28 .. 2021-06-28) to the end of the existing data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44344, 0, 0, 0),
    @(44345, 0, 0, 0),
    @(44346, 0, 0, 0),
    @(44347, 0, 0, 0),
    @(44348, 0, 0, 0),
    @(44349, 0, 0, 0),
    @(44350, 0, 0, 0),
    @(44351, 0, 0, 0),
    @(44352, 0, 0, 0),
    @(44353, 0, 0, 0),
    @(44354, 0, 0, 0),
    @(44355, 1, 1, 37.46721618583739),
    @(44356, 0, 1, 37.46721618583739),
    @(44357, 0, 1, 37.46721618583739),
    @(44358, 0, 1, 37.46721618583739),
    @(44359, 0, 1, 37.46721618583739),
    @(44360, 0, 1, 37.46721618583739),
    @(44361, 0, 1, 37.46721618583739),
    @(44362, 0, 0, 0),
    @(44363, 0, 0, 0),
    @(44364, 0, 0, 0),
    @(44365, 0, 0, 0),
    @(44366, 0, 0, 0),
    @(44367, 0, 0, 0),
    @(44368, 0, 0, 0),
    @(44369, 0, 0, 0),
    @(44370, 0, 0, 0),
    @(44371, 0, 0, 0),
    @(44372, 0, 0, 0),
    @(44373, 0, 0, 0),
    @(44374, 0, 0, 0),
    @(44375, 0, 0, 0)
)

$lastRow = 269
$startRow = $lastRow + 1
$row = $startRow

foreach ($entry in $newRows) {
    # Column A keeps the same date style as the row above it (bold, bordered,
    # centered, custom date number format) - copy it, then overwrite the value.
    $ws.Range("A$lastRow").Copy($ws.Range("A$row"))
    $ws.Range("A$row").Value = $entry[0]

    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $ws.Range("D$row").Value = $entry[3]

    $row = $row + 1
}

Write-Output ("Updated range: " + $ws.UsedRange.Address())
